$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.748.31"
$ws.Range("D3").Value = "2.624.03"
$ws.Range("E3").Value = "  -0.59%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "596.13"
$ws.Range("E5").Value = "  -1.34%  "
$ws.Range("D6").Value = "150.74"
$ws.Range("E6").Value = "  +2.91%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "0.590"
$ws.Range("E8").Value = "  +0.27%  "
$ws.Range("E9").Value = "  +0.62%  "
$ws.Range("E10").Value = "  +1.79%  "
$ws.Range("E11").Value = "  +3.52%  "
$ws.Range("E12").Value = "  -1.06%  "
$ws.Range("E13").Value = "  +0.28%  "
$ws.Range("D14").Value = "3.092.40"
$ws.Range("E14").Value = "  -0.65%  "
$ws.Range("D15").Value = "63.558.86"
$ws.Range("E15").Value = "  +0.34%  "
$ws.Range("D16").Value = "0.0000154"
$ws.Range("E16").Value = "  +4.94%  "
$ws.Range("D17").Value = "2.621.88"
$ws.Range("E17").Value = "  -0.62%  "
$ws.Range("D18").Value = "12.32"
$ws.Range("E18").Value = "  +6.95%  "
$ws.Range("E19").Value = "  +1.81%  "
$ws.Range("D20").Value = "347.76"
$ws.Range("E20").Value = "  +1.07%  "
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").Value = "5.72"
$ws.Range("E23").Value = "  +2.77%  "
$ws.Range("D24").Value = "66.41"
$ws.Range("E24").Value = "  -0.31%  "
$ws.Range("D25").Value = "1.72"
$ws.Range("E25").Value = "  +12.17%  "
$ws.Range("D26").Value = "9.26"
$ws.Range("E26").Value = "  +2.03%  "
$ws.Range("E27").Value = "  -0.86%  "
$ws.Range("D28").Value = "567.77"
$ws.Range("E28").Value = "  -3.38%  "
$ws.Range("D29").Value = "8.26"
$ws.Range("E29").Value = "  +3.47%  "
$ws.Range("E30").Value = "  -0.42%  "
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("E32").Value = "  -0.75%  "
$ws.Range("D33").Value = "0.0₃0848"
$ws.Range("E33").Value = "  +2.10%  "
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("D35").Value = "5.23"
$ws.Range("E35").Value = "  +0.22%  "
$ws.Range("D36").Value = "168.94"
$ws.Range("E36").Value = "  +1.34%  "
$ws.Range("E37").Value = "  +0.57%  "
$ws.Range("D39").Value = "1.94"
$ws.Range("E39").Value = "  -0.94%  "
$ws.Range("E40").Value = "  +1.58%  "
$ws.Range("D42").Value = "166.64"
$ws.Range("E42").Value = "  -0.85%  "
$ws.Range("D43").Value = "39.93"
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("E44").Value = "  +3.49%  "
$ws.Range("D45").Value = "0.0595"
$ws.Range("E45").Value = "  +4.44%  "
$ws.Range("D46").Value = "21.64"
$ws.Range("E46").Value = "  -2.51%  "
$ws.Range("E47").Value = "  +0.26%  "
$ws.Range("E48").Value = "  +1.42%  "
$ws.Range("D49").Value = "2.00"
$ws.Range("E49").Value = "  +3.94%  "
$ws.Range("E50").Value = "  +0.22%  "
$ws.Range("D51").Value = "19.44"
$ws.Range("E51").Value = "  +3.46%  "
